# Pinout sheet update:
#  - Add an external TAS "Replay Toggle" button wired to INT0 (row 4, cols I/J/K)
#  - Add a status LED wired to VDD/port info (row 11, cols B/C)
#  - Slightly widen the two outer "Name" columns (A and K)
#  - Leave the active selection on A12

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New pin: INT0 used as an external "Replay Toggle" (TAS start) button input
$ws.Range("I4").Value = "INT0"
$ws.Range("J4").Value = "I"
$ws.Range("K4").Value = "Replay Toggle"

# New pin: status LED output
$ws.Range("B11").Value = "O"
$ws.Range("C11").Value = "STAT_LED"

# Widen column A and column K slightly to fit the new labels
$ws.Columns.Item(1).ColumnWidth = 15.71
$ws.Columns.Item(11).ColumnWidth = 15.71

# Move the active selection
$ws.Range("A12").Select() | Out-Null
